$d = $word.ActiveDocument

# Locate the paragraph containing "Explanation of the schema " and expand
# the range to cover the whole paragraph (including its paragraph mark).
$anchor = $d.Content
$found = $anchor.Find.Execute("Explanation of the schema", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor paragraph 'Explanation of the schema'"
}
$anchor.Expand(4) | Out-Null   # wdParagraph

# --- Step 1: give the existing run an explicit (empty) <w:rPr/> ------------
# Target just the run's text (exclude the trailing paragraph mark) and
# re-insert the identical text wrapped with an empty run-properties element,
# so the run gains a <w:rPr/> without altering its content.
$runRange = $d.Range($anchor.Start, $anchor.End - 1)
$existingText = $runRange.Text
$runXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
      '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
          '<w:body><w:p><w:r><w:rPr/><w:t xml:space="preserve">' + $existingText + '</w:t></w:r></w:p></w:body>' + `
        '</w:document>' + `
      '</pkg:xmlData>' + `
    '</pkg:part>' + `
  '</pkg:package>'
$runRange.InsertXML($runXml)

# --- Step 2: append four new "Normal" style paragraphs right after it ------
$anchor = $d.Content
$found = $anchor.Find.Execute("Explanation of the schema", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not re-find anchor paragraph 'Explanation of the schema'"
}
$anchor.Expand(4) | Out-Null   # wdParagraph
$insertPoint = $d.Range($anchor.End, $anchor.End)

$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
      '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
          '<w:body>' + `
            '<w:p><w:pPr><w:pStyle w:val="Normal"/></w:pPr><w:r><w:rPr/><w:t>Purchasing licence</w:t></w:r></w:p>' + `
            '<w:p><w:pPr><w:pStyle w:val="Normal"/></w:pPr><w:r><w:rPr/><w:t>Pricing of licences</w:t></w:r></w:p>' + `
            '<w:p><w:pPr><w:pStyle w:val="Normal"/></w:pPr><w:r><w:rPr/><w:t>Customers and country</w:t></w:r></w:p>' + `
            '<w:p><w:pPr><w:pStyle w:val="Normal"/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">Naming convention </w:t></w:r></w:p>' + `
          '</w:body>' + `
        '</w:document>' + `
      '</pkg:xmlData>' + `
    '</pkg:part>' + `
  '</pkg:package>'
$insertPoint.InsertXML($newParaXml)

Write-Output "Inserted 4 paragraphs; document now has $($d.Paragraphs.Count) paragraphs."
